$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Fix Email HOLD/SELL Recommendations (Cache Invalidation)"
$ws.Range("C8").Value = "Email Scanner"
$ws.Range("D8").Value = "Bug critico: le raccomandazioni HOLD/SELL nell'email giornaliera erano basate su dati cached obsoleti (TICKER_CACHE). Questo causava 3 errori su 9 posizioni: RACE.MI e MUV2.DE mostravano HOLD invece di SELL (strategie chiuse l'11-12/02), MMM mostrava SELL invece di HOLD (strategia aperta dall'11/02). Secondo bug: ticker del portafoglio non presenti in tickers.js non venivano mai scansionati, risultando sempre in SELL."
$ws.Range("E8").Value = "1) Creazione test standalone (test_email_standalone.py) con caricamento portfolio da Firebase. 2) Analisi indipendente di tutti i 9 ticker del portafoglio per confermare le discrepanze. 3) Invalidazione mirata della TICKER_CACHE per i soli ticker del portafoglio prima dello scan. 4) Merge automatico dei ticker del portafoglio nella lista di scan."
$ws.Range("F8").Value = "backend/scanner.py, test_email_standalone.py"
$ws.Range("G8").Value = "Done"
$ws.Range("H8").Value = "Critical"
$ws.Range("I8").Value = "NO"
$ws.Range("J8").Value = "Data Consistency, Email Accuracy"
$ws.Range("K8").Value = "SI"
$ws.Range("L8").Value = "1) test_email_standalone.py - confronto raccomandazioni con email reale. 2) Verifica post-deploy: email delle 16:30 deve riflettere stato reale delle strategie."
$ws.Range("M8").Value = "Test standalone conferma: 4 HOLD (ITW, JNJ, TTE.PA, PANW), 2 SELL (RACE.MI, MUV2.DE). Tutte le discrepanze con l'email risolte."
$ws.Range("N8").Value = "Deploy ID: PENDING (15/02/2026)"
